function Set-ParagraphInnerXml {
    param($Document, [string]$AnchorText, [string]$InnerXml)

    $r = $Document.Content
    $r.Find.ClearFormatting()
    $found = $r.Find.Execute($AnchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $AnchorText"
    }
    $r.Expand(4) | Out-Null
    $target = $Document.Range($r.Start, $r.End - 1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

$d = $word.ActiveDocument

# 1) "Date:" paragraph -- add the "March 17, 2023" run after the trailing space run.
Set-ParagraphInnerXml $d "Date: " (
    '<w:r w:rsidRPr="00B556F4"><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:szCs w:val="24"/></w:rPr><w:t>Date:</w:t></w:r>' +
    '<w:r w:rsidRPr="00B556F4"><w:rPr><w:rFonts w:cs="Times New Roman"/><w:i/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:iCs/><w:szCs w:val="24"/></w:rPr><w:t>March 17, 2023</w:t></w:r>'
)

# 2) "This will generate a forecast report..." paragraph -- several edits:
#    - ".html file" -> "Word document (.docx)" (split into 3 runs)
#    - the quoted "figures" folder text -> italicized figures run (split into 4 runs)
#    - "fit new models " / "before generating" -> "fit new " / "models before generating"
Set-ParagraphInnerXml $d "This will generate a forecast report" (
    '<w:r><w:t xml:space="preserve">This will generate a forecast report for the desired year and house the report and figures within directories for that year. The report will be output as a </w:t></w:r>' +
    '<w:r><w:t>Word document (.docx)</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> located in </w:t></w:r>' +
    '<w:r w:rsidR="00696086"><w:t>a folder for that year within the</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidRPr="0064593E"><w:rPr><w:i/><w:iCs/></w:rPr><w:t>report output</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> folder. </w:t></w:r>' +
    '<w:r w:rsidR="00696086"><w:t xml:space="preserve">It will also produce figures for model observed and predicted values within the </w:t></w:r>' +
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>figures</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">folder. </w:t></w:r>' +
    '<w:r w:rsidR="00696086"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Note: </w:t></w:r>' +
    '<w:r w:rsidR="006D6E6B"><w:t xml:space="preserve">This process will take a few minutes to fit new </w:t></w:r>' +
    '<w:r w:rsidR="006D6E6B"><w:lastRenderedPageBreak/><w:t xml:space="preserve">models before generating the report. </w:t></w:r>' +
    '<w:r w:rsidR="006D6E6B" w:rsidRPr="006D6E6B"><w:t>Additionally, i</w:t></w:r>' +
    '<w:r w:rsidR="00696086"><w:t xml:space="preserve">f this is the first time running the report on a </w:t></w:r>' +
    '<w:r w:rsidR="006D6E6B"><w:t>computer, it will take a few minutes to download and install the necessary R packages.</w:t></w:r>'
)

# 3) "...delete the model objects..." paragraph -- append " To do so:" run.
Set-ParagraphInnerXml $d "Save and close out the spreadsheet" (
    '<w:r><w:t xml:space="preserve">Save and close out the spreadsheet, then run the report. </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Note: </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">If the report has already been run for </w:t></w:r>' +
    '<w:r w:rsidR="001D3551"><w:t xml:space="preserve">the </w:t></w:r>' +
    '<w:r><w:t>year, you will need to delete the model objects that were created in order to fit models with the new data.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> To do so:</w:t></w:r>'
)

# 4) "To do so, navigate to..." paragraph -- split leading "To do so, navigate to " into "N" + "avigate to ".
Set-ParagraphInnerXml $d "To do so, navigate to" (
    '<w:r><w:t>N</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">avigate to </w:t></w:r>' +
    '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>R/model objects/</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and select the folder with the prediction year.</w:t></w:r>'
)

Write-Output "All edits applied."
